$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '53.200.02'
$ws.Range('E2').Value = '  -8.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.381.28'
$ws.Range('E3').Value = '  -11.77%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '457.63'
$ws.Range('E5').Value = '  -8.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.26'
$ws.Range('E6').Value = '  -6.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.479'
$ws.Range('E8').Value = '  -9.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.398.18'
$ws.Range('E9').Value = '  -11.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0941'
$ws.Range('E10').Value = '  -7.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.23'
$ws.Range('E11').Value = '  -12.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.311'
$ws.Range('E12').Value = '  -8.89%  '
$ws.Range('E13').Value = '  -4.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.805.40'
$ws.Range('E14').Value = '  -11.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '53.361.00'
$ws.Range('E15').Value = '  -8.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.51'
$ws.Range('E16').Value = '  -8.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000129'
$ws.Range('E17').Value = '  -3.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.398.42'
$ws.Range('E18').Value = '  -10.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.13'
$ws.Range('E19').Value = '  -11.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '304.64'
$ws.Range('E20').Value = '  -9.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.24'
$ws.Range('E21').Value = '  -14.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.994'
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.67'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.29'
$ws.Range('E24').Value = '  -13.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '55.69'
$ws.Range('E25').Value = '  -10.71%  '
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.380'
$ws.Range('E27').Value = '  -9.19%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.518.80'
$ws.Range('E28').Value = '  -10.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.152'
$ws.Range('E29').Value = '  -10.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.05'
$ws.Range('E30').Value = '  -4.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0₃0707'
$ws.Range('E32').Value = '  -13.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '145.40'
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.52'
$ws.Range('E34').Value = '  -7.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.41'
$ws.Range('E35').Value = '  -11.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.94'
$ws.Range('E36').Value = '  -6.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.48'
$ws.Range('E37').Value = '  -15.31%  '
$ws.Range('E38').Value = '  -6.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.783'
$ws.Range('E39').Value = '  -14.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.996'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.85'
$ws.Range('E41').Value = '  -8.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.590'
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0519'
$ws.Range('E43').Value = '  -5.62%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.23'
$ws.Range('E44').Value = '  -6.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.16'
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.22'
$ws.Range('E46').Value = '  -10.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.929.59'
$ws.Range('E47').Value = '  -10.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0214'
$ws.Range('E48').Value = '  -4.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0858'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.15'
$ws.Range('E50').Value = '  -9.20%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.36'
$ws.Range('E51').Value = '  -11.81%  '
